# Updates cryptos price list (columns D = Price, E = Volume(1h))
# for the rows whose values moved since the last snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.246.33"
$ws.Range("E2").Value2 = "  -0.34%  "
$ws.Range("D3").Value = "'1.812.17"
$ws.Range("E3").Value2 = "  +1.45%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value2 = "  -0.29%  "
$ws.Range("D5").Value = "'225.25"
$ws.Range("E5").Value2 = "  +0.23%  "
$ws.Range("D6").Value = "'0.556"
$ws.Range("E6").Value2 = "  +0.80%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value2 = "  -0.22%  "
$ws.Range("D8").Value = "'32.23"
$ws.Range("E8").Value2 = "  -1.53%  "
$ws.Range("D9").Value = "'0.295"
$ws.Range("E9").Value2 = "  +4.73%  "
$ws.Range("D10").Value = "'0.0728"
$ws.Range("E10").Value2 = "  +10.22%  "
$ws.Range("D11").Value = "'0.0929"
$ws.Range("E11").Value2 = "  -0.33%  "
$ws.Range("D12").Value = "'2.073.12"
$ws.Range("E12").Value2 = "  +1.45%  "
$ws.Range("D13").Value = "'1.822.91"
$ws.Range("E13").Value2 = "  +2.23%  "
$ws.Range("D14").Value = "'10.98"
$ws.Range("E14").Value2 = "  -0.20%  "
$ws.Range("D15").Value = "'0.640"
$ws.Range("E15").Value2 = "  +0.86%  "
$ws.Range("D16").Value = "'34.210.68"
$ws.Range("E16").Value2 = "  -0.47%  "
$ws.Range("E17").Value2 = "  +1.66%  "
$ws.Range("D18").Value = "'69.61"
$ws.Range("E18").Value2 = "  +0.77%  "
$ws.Range("D19").Value = "'249.45"
$ws.Range("E19").Value2 = "  -2.17%  "
$ws.Range("D20").Value = "'0.0₃0799"
$ws.Range("E20").Value2 = "  +7.31%  "
$ws.Range("D21").Value = "'11.01"
$ws.Range("E21").Value2 = "  +6.21%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value2 = "  -0.18%  "
$ws.Range("D23").Value = "'4.25"
$ws.Range("E23").Value2 = "  +1.04%  "
$ws.Range("E24").Value2 = "  +0.57%  "
$ws.Range("D25").Value = "'160.25"
$ws.Range("E25").Value2 = "  +1.73%  "
$ws.Range("D26").Value = "'16.72"
$ws.Range("E26").Value2 = "  +1.75%  "
$ws.Range("D27").Value = "'7.20"
$ws.Range("E27").Value2 = "  +2.75%  "
$ws.Range("E28").Value2 = "  +0.71%  "
$ws.Range("D29").Value = "'0.997"
$ws.Range("E29").Value2 = "  -0.42%  "
$ws.Range("D30").Value = "'0.0533"
$ws.Range("E30").Value2 = "  +3.82%  "
$ws.Range("D31").Value = "'3.78"
$ws.Range("E31").Value2 = "  +0.10%  "
$ws.Range("D32").Value = "'1.22"
$ws.Range("E32").Value2 = "  +2.27%  "
$ws.Range("E33").Value2 = "  -0.15%  "
$ws.Range("E34").Value2 = "  -0.39%  "
$ws.Range("D35").Value = "'1.432.46"
$ws.Range("E35").Value2 = "  -1.42%  "
$ws.Range("E36").Value2 = "  +1.15%  "
$ws.Range("D37").Value = "'0.640"
$ws.Range("E37").Value2 = "  +1.64%  "
$ws.Range("E38").Value2 = "  +0.54%  "
$ws.Range("D39").Value = "'0.960"
$ws.Range("E39").Value2 = "  +7.94%  "
$ws.Range("D40").Value = "'81.22"
$ws.Range("E40").Value2 = "  -2.58%  "
$ws.Range("E41").Value2 = "  -3.65%  "
$ws.Range("E43").Value2 = "  +4.06%  "
$ws.Range("E44").Value2 = "  +1.95%  "
$ws.Range("E45").Value2 = "  -1.58%  "
$ws.Range("D46").Value = "'1.970.62"
$ws.Range("D47").Value = "'1.05"
$ws.Range("E47").Value2 = "  -0.71%  "
$ws.Range("D48").Value = "'106.70"
$ws.Range("E48").Value2 = "  +7.71%  "
$ws.Range("D49").Value = "'11.99"
$ws.Range("E49").Value2 = "  -2.04%  "
$ws.Range("E50").Value2 = "  -0.48%  "
$ws.Range("E51").Value2 = "  +6.64%  "
